# Auto-generated script to update cryptos price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextValue $ws.Cells.Item(2, 4) '28.653.08'
Set-TextValue $ws.Cells.Item(2, 5) '  +2.21%  '
Set-TextValue $ws.Cells.Item(3, 4) '1.800.34'
Set-TextValue $ws.Cells.Item(3, 5) '  -0.16%  '
Set-TextValue $ws.Cells.Item(4, 4) '1.001'
Set-TextValue $ws.Cells.Item(4, 5) '  +0.23%  '
Set-TextValue $ws.Cells.Item(5, 4) '313.66'
Set-TextValue $ws.Cells.Item(5, 5) '  -0.61%  '
Set-TextValue $ws.Cells.Item(6, 5) '  +0.27%  '
Set-TextValue $ws.Cells.Item(7, 4) '0.5397'
Set-TextValue $ws.Cells.Item(7, 5) '  -0.96%  '
Set-TextValue $ws.Cells.Item(8, 4) '0.3779'
Set-TextValue $ws.Cells.Item(8, 5) '  -0.42%  '
Set-TextValue $ws.Cells.Item(9, 4) '0.07534'
Set-TextValue $ws.Cells.Item(9, 5) '  -0.64%  '
Set-TextValue $ws.Cells.Item(10, 4) '42.59'
Set-TextValue $ws.Cells.Item(10, 5) '  -1.11%  '
Set-TextValue $ws.Cells.Item(11, 4) '1.116'
Set-TextValue $ws.Cells.Item(11, 5) '  -1.46%  '
Set-TextValue $ws.Cells.Item(12, 5) '  +0.33%  '
Set-TextValue $ws.Cells.Item(13, 4) '20.96'
Set-TextValue $ws.Cells.Item(13, 5) '  -0.83%  '
Set-TextValue $ws.Cells.Item(14, 4) '6.181'
Set-TextValue $ws.Cells.Item(14, 5) '  -0.80%  '
Set-TextValue $ws.Cells.Item(15, 4) '7.426'
Set-TextValue $ws.Cells.Item(15, 5) '  +4.03%  '
Set-TextValue $ws.Cells.Item(16, 4) '1.796.39'
Set-TextValue $ws.Cells.Item(16, 5) '  -0.18%  '
Set-TextValue $ws.Cells.Item(17, 4) '90.60'
Set-TextValue $ws.Cells.Item(17, 5) '  -1.66%  '
Set-TextValue $ws.Cells.Item(18, 4) '0.00001066'
Set-TextValue $ws.Cells.Item(18, 5) '  -1.20%  '
Set-TextValue $ws.Cells.Item(19, 4) '0.06446'
Set-TextValue $ws.Cells.Item(19, 5) '  -0.88%  '
Set-TextValue $ws.Cells.Item(20, 5) '  +0.20%  '
Set-TextValue $ws.Cells.Item(21, 4) '17.27'
Set-TextValue $ws.Cells.Item(21, 5) '  +0.71%  '
Set-TextValue $ws.Cells.Item(22, 4) '5.938'
Set-TextValue $ws.Cells.Item(22, 5) '  -0.73%  '
Set-TextValue $ws.Cells.Item(23, 4) '28.647.84'
Set-TextValue $ws.Cells.Item(23, 5) '  +2.12%  '
Set-TextValue $ws.Cells.Item(24, 4) '11.18'
Set-TextValue $ws.Cells.Item(24, 5) '  -0.53%  '
Set-TextValue $ws.Cells.Item(25, 4) '2.102'
Set-TextValue $ws.Cells.Item(26, 4) '160.67'
Set-TextValue $ws.Cells.Item(26, 5) '  +2.70%  '
Set-TextValue $ws.Cells.Item(27, 4) '20.50'
Set-TextValue $ws.Cells.Item(27, 5) '  -0.45%  '
Set-TextValue $ws.Cells.Item(28, 4) '2.378'
Set-TextValue $ws.Cells.Item(28, 5) '  -0.50%  '
Set-TextValue $ws.Cells.Item(29, 4) '2.001.30'
Set-TextValue $ws.Cells.Item(29, 5) '  -0.37%  '
Set-TextValue $ws.Cells.Item(30, 4) '123.54'
Set-TextValue $ws.Cells.Item(30, 5) '  +0.48%  '
Set-TextValue $ws.Cells.Item(31, 4) '1.109'
Set-TextValue $ws.Cells.Item(31, 5) '  -3.75%  '
Set-TextValue $ws.Cells.Item(32, 4) '0.1044'
Set-TextValue $ws.Cells.Item(32, 5) '  +0.87%  '
Set-TextValue $ws.Cells.Item(33, 4) '5.673'
Set-TextValue $ws.Cells.Item(33, 5) '  -1.44%  '
Set-TextValue $ws.Cells.Item(34, 4) '3.694'
Set-TextValue $ws.Cells.Item(34, 5) '  +2.59%  '
Set-TextValue $ws.Cells.Item(35, 4) '0.2260'
Set-TextValue $ws.Cells.Item(35, 5) '  +6.77%  '
Set-TextValue $ws.Cells.Item(36, 4) '0.06504'
Set-TextValue $ws.Cells.Item(36, 5) '  +7.27%  '
Set-TextValue $ws.Cells.Item(37, 4) '8.889'
Set-TextValue $ws.Cells.Item(37, 5) '  +2.63%  '
Set-TextValue $ws.Cells.Item(38, 4) '0.02322'
Set-TextValue $ws.Cells.Item(38, 5) '  +0.98%  '
Set-TextValue $ws.Cells.Item(39, 4) '5.047'
Set-TextValue $ws.Cells.Item(39, 5) '  +0.36%  '
Set-TextValue $ws.Cells.Item(40, 4) '11.36'
Set-TextValue $ws.Cells.Item(40, 5) '  -1.27%  '
Set-TextValue $ws.Cells.Item(41, 4) '1.212'
Set-TextValue $ws.Cells.Item(41, 5) '  +5.24%  '
Set-TextValue $ws.Cells.Item(42, 4) '0.6254'
Set-TextValue $ws.Cells.Item(42, 5) '  -0.81%  '
Set-TextValue $ws.Cells.Item(43, 4) '1.000'
Set-TextValue $ws.Cells.Item(43, 5) '  +0.19%  '
Set-TextValue $ws.Cells.Item(44, 5) '  -0.22%  '
Set-TextValue $ws.Cells.Item(45, 4) '13.42'
Set-TextValue $ws.Cells.Item(45, 5) '  -0.10%  '
Set-TextValue $ws.Cells.Item(46, 4) '0.5880'
Set-TextValue $ws.Cells.Item(46, 5) '  -0.79%  '
Set-TextValue $ws.Cells.Item(47, 4) '3.665'
Set-TextValue $ws.Cells.Item(47, 5) '  -0.11%  '
Set-TextValue $ws.Cells.Item(48, 4) '126.41'
Set-TextValue $ws.Cells.Item(48, 5) '  +3.42%  '
Set-TextValue $ws.Cells.Item(49, 4) '1.961'
Set-TextValue $ws.Cells.Item(49, 5) '  +1.79%  '
Set-TextValue $ws.Cells.Item(50, 4) '1.160'
Set-TextValue $ws.Cells.Item(50, 5) '  +2.18%  '
Set-TextValue $ws.Cells.Item(51, 4) '0.06897'
Set-TextValue $ws.Cells.Item(51, 5) '  +1.61%  '
